$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the "Total" row (currently row 16) down to row 19, leaving rows
# 15-17 free for the new entries (row 18 stays blank, just like row 9/12
# elsewhere in the sheet) and carrying the row's formatting (border,
# height, styles) along with it.
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()

# Row 15 is left blank except for the date-column styling, matching the
# existing blank separator rows (e.g. row 12).
$ws.Range("B15").NumberFormat = "m/d/yy"
$ws.Range("B15").HorizontalAlignment = -4108

# Row 16 - new 2008 data entry
$ws.Range("A16").Value = "Added all the 2008 data to a modified database"
$ws.Range("B16").Value = 39709
$ws.Range("B16").NumberFormat = "m/d/yy"
$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("C16").Value = 2
$ws.Range("C16").HorizontalAlignment = -4108

# Row 17 - new 2008 data entry
$ws.Range("A17").Value = "Modified web pages as suggested my Marting"
$ws.Range("B17").Value = 39710
$ws.Range("B17").NumberFormat = "m/d/yy"
$ws.Range("B17").HorizontalAlignment = -4108
$ws.Range("C17").Value = 1
$ws.Range("C17").HorizontalAlignment = -4108

# Extend the Total formula to cover the new rows.
$ws.Range("C19").Formula = "=SUM(C2:C18)"

# Update the view: selection moves to the new last data cell.
$ws.Range("C18").Select() | Out-Null
